$d = $word.ActiveDocument

# Locate the full placeholder phrase "(${jmlhari} hari sejak di terbitkan)"
# which today is spread across six separate runs: "(", "$", "{", "jmlhari",
# "}", " hari sejak di terbitkan)". We need to collapse it down to just
# "${jmlhari}" while keeping that as four individual runs (one per token:
# "$", "{", "jmlhari", "}") so that each run keeps the same run-level
# formatting (Arial Narrow, bold, color 000000, lang nb-NO) as before.
$rng = $d.Content
$found = $rng.Find.Execute('(${jmlhari} hari sejak di terbitkan)')

if ($found) {
    $runPr = '<w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:val="nb-NO"/></w:rPr>'

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p>' +
           '<w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/>' + $runPr + '</w:pPr>' +
           '<w:r>' + $runPr + '<w:t>$</w:t></w:r>' +
           '<w:r>' + $runPr + '<w:t>{</w:t></w:r>' +
           '<w:r>' + $runPr + '<w:t>jmlhari</w:t></w:r>' +
           '<w:r>' + $runPr + '<w:t>}</w:t></w:r>' +
           '</w:p>' +
           '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
}
